# Update three data rows that previously duplicated another row so that
# they contain a brand-new (distinct) person record. This matches the
# "repeated data" sorting test-case scenario described in the commit
# message (FMT6 / Test_FMT4_C - sorting with repeated column values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was a duplicate of row 3 (Doug Derrick / Timepath Inc.). Replace it
# with a new "Lara Palmer" record at the same company.
$ws.Range("A2").Value = "Lara"
$ws.Range("B2").Value = "Palmer"
$ws.Range("C2").Value = "Timepath Inc."
$ws.Range("D2").Value = "Programmer"
$ws.Range("E2").Value = "87 Orange Street"
$ws.Range("F2").Value = "lpalmer@timepath.co.uk"
$ws.Range("G2").Value = 40731653845

# Row 6 was a duplicate of row 7 (Michael Robertson / MediCare). Replace it
# with a new "Jane Dorsey" record at the same company.
$ws.Range("A6").Value = "Jane"
$ws.Range("B6").Value = "Dorsey"
$ws.Range("C6").Value = "MediCare"
$ws.Range("D6").Value = "Medical Engineer"
$ws.Range("E6").Value = "11 Crown Street"
$ws.Range("F6").Value = "jdorsey@mc.com"
$ws.Range("G6").Value = 40791345621

# Row 9 was a duplicate of row 10 (Jessie Marlowe / Aperture Inc.). Replace
# it with a new "Michelle Norton" record at the same company/role.
$ws.Range("A9").Value = "Michelle"
$ws.Range("B9").Value = "Norton"
$ws.Range("C9").Value = "Aperture Inc."
$ws.Range("D9").Value = "Scientist"
$ws.Range("E9").Value = "13 White Rabbit Street"
$ws.Range("F9").Value = "mnorton@aperture.us"
$ws.Range("G9").Value = 40731254562

$wb.Save()
